# Updates cryptocurrency price/volume data to the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '44.017.43'
Set-TextValue 'E2' '  +0.09%  '

Set-TextValue 'D3' '2.239.45'
Set-TextValue 'E3' '  -0.31%  '

Set-TextValue 'E4' '  +0.22%  '

Set-TextValue 'D5' '306.55'
Set-TextValue 'E5' '  -3.83%  '

Set-TextValue 'D6' '95.14'
Set-TextValue 'E6' '  -5.09%  '

Set-TextValue 'D7' '0.571'
Set-TextValue 'E7' '  -0.40%  '

Set-TextValue 'E8' '  +0.25%  '

Set-TextValue 'D9' '0.523'
Set-TextValue 'E9' '  -3.75%  '

Set-TextValue 'D10' '34.75'
Set-TextValue 'E10' '  -4.94%  '

Set-TextValue 'E11' '  -1.75%  '

Set-TextValue 'D12' '7.21'
Set-TextValue 'E12' '  -3.61%  '

Set-TextValue 'E13' '  -0.78%  '

Set-TextValue 'D14' '2.581.27'
Set-TextValue 'E14' '  -0.17%  '

Set-TextValue 'D15' '2.235.59'
Set-TextValue 'E15' '  +0.45%  '

Set-TextValue 'D16' '0.824'
Set-TextValue 'E16' '  -2.74%  '

Set-TextValue 'D17' '13.57'
Set-TextValue 'E17' '  -4.38%  '

Set-TextValue 'D18' '43.920.89'
Set-TextValue 'E18' '  +0.13%  '

Set-TextValue 'D19' '0.0₃0964'
Set-TextValue 'E19' '  -0.72%  '

Set-TextValue 'D20' '12.06'
Set-TextValue 'E20' '  -10.29%  '

Set-TextValue 'D21' '6.29'
Set-TextValue 'E21' '  -2.11%  '

Set-TextValue 'D22' '65.26'
Set-TextValue 'E22' '  +0.01%  '

Set-TextValue 'D23' '236.95'
Set-TextValue 'E23' '  +1.03%  '

Set-TextValue 'D24' '2.93'
Set-TextValue 'E24' '  -4.71%  '

Set-TextValue 'E25' '  -4.54%  '

Set-TextValue 'E26' '  +0.22%  '

Set-TextValue 'D27' '9.93'
Set-TextValue 'E27' '  -4.90%  '

Set-TextValue 'B28' 'InjectiveProtocol'
Set-TextValue 'C28' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D28' '37.87'
Set-TextValue 'E28' '  -1.23%  '

Set-TextValue 'B29' 'Toncoin'
Set-TextValue 'C29' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D29' '2.17'
Set-TextValue 'E29' '  -1.65%  '

Set-TextValue 'D30' '6.03'
Set-TextValue 'E30' '  -0.22%  '

Set-TextValue 'E31' '  -0.77%  '

Set-TextValue 'D32' '152.90'
Set-TextValue 'E32' '  -3.67%  '

Set-TextValue 'D33' '0.0802'
Set-TextValue 'E33' '  -5.07%  '

Set-TextValue 'D34' '3.27'
Set-TextValue 'E34' '  +2.56%  '

Set-TextValue 'E35' '  -2.83%  '

Set-TextValue 'E36' '  -4.15%  '

Set-TextValue 'E37' '  +0.84%  '

Set-TextValue 'D38' '1.78'
Set-TextValue 'E38' '  -7.78%  '

Set-TextValue 'D39' '15.08'
Set-TextValue 'E39' '  -6.18%  '

Set-TextValue 'D40' '3.83'
Set-TextValue 'E40' '  -7.12%  '

Set-TextValue 'D41' '3.33'
Set-TextValue 'E41' '  -8.75%  '

Set-TextValue 'E42' '  -3.96%  '

Set-TextValue 'D44' '1.725.17'
Set-TextValue 'E44' '  -2.00%  '

Set-TextValue 'D45' '85.16'
Set-TextValue 'E45' '  +4.86%  '

Set-TextValue 'D46' '0.188'
Set-TextValue 'E46' '  -3.64%  '

Set-TextValue 'D47' '100.06'
Set-TextValue 'E47' '  -2.84%  '

Set-TextValue 'D48' '4.94'
Set-TextValue 'E48' '  -3.70%  '

Set-TextValue 'D49' '68.89'
Set-TextValue 'E49' '  -7.48%  '

Set-TextValue 'D50' '8.06'
Set-TextValue 'E50' '  -2.70%  '

Set-TextValue 'B51' 'EnergySwap'
Set-TextValue 'C51' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D51' '14.33'
Set-TextValue 'E51' '  -0.40%  '
